$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "13:00-zeit"
$ws.Range("A11").Value = "22.12.2023"
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat

$ws.Range("A11").Select()
